$wb = $excel.ActiveWorkbook

foreach ($sheet in $wb.Worksheets) {
    $sheet.Name = $sheet.Name -replace '^(data_(?:CCM|RAM)) code_', '$1-code_'
}
